# Correct the "Банка" (Bank) column: the rows that list FULM Shtedilnica
# products were mistakenly carrying the "Стопанска Банка АД Скопје" label.
# Update them to the correct bank name. Rows 6 and 9 already show a
# different product and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bankName = "ФУЛМ Штедилница д.о.о. - Скопје"
$rows = @(2, 3, 4, 5, 7, 8, 10, 11, 12, 13, 14, 15)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 1).Value = $bankName
}

# Restore the view to where the author left it: scrolled down so row 12
# is at the top, with the blank rows below the table selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
$ws.Rows("16:1048576").Select()
